$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-04-08 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-09 Wednesday", 2)

# Update the answer table cell by cell
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "68-54=14"
$t.Cell(1, 2).Range.Text = "17+66=83"
$t.Cell(1, 3).Range.Text = "31+42=73"
$t.Cell(1, 4).Range.Text = "81-37=44"
$t.Cell(1, 5).Range.Text = "84-33=51"
$t.Cell(2, 1).Range.Text = "2+81=83"
$t.Cell(2, 2).Range.Text = "65-61=4"
$t.Cell(2, 3).Range.Text = "75-56=19"
$t.Cell(2, 4).Range.Text = "22+32=54"
$t.Cell(2, 5).Range.Text = "60-23=37"
$t.Cell(3, 1).Range.Text = "38+54=92"
$t.Cell(3, 2).Range.Text = "51-1=50"
$t.Cell(3, 3).Range.Text = "15+52=67"
$t.Cell(3, 4).Range.Text = "13-11=2"
$t.Cell(3, 5).Range.Text = "41+5=46"
$t.Cell(4, 1).Range.Text = "71-39=32"
$t.Cell(4, 2).Range.Text = "93-83=10"
$t.Cell(4, 3).Range.Text = "66-48=18"
$t.Cell(4, 4).Range.Text = "30+57=87"
$t.Cell(4, 5).Range.Text = "41+39=80"
$t.Cell(5, 1).Range.Text = "90-2=88"
$t.Cell(5, 2).Range.Text = "94-35=59"
$t.Cell(5, 3).Range.Text = "21+43=64"
$t.Cell(5, 4).Range.Text = "97-15=82"
$t.Cell(5, 5).Range.Text = "51+35=86"
$t.Cell(6, 1).Range.Text = "96-3=93"
$t.Cell(6, 2).Range.Text = "48+46=94"
$t.Cell(6, 3).Range.Text = "22+68=90"
$t.Cell(6, 4).Range.Text = "55+12=67"
$t.Cell(6, 5).Range.Text = "75-8=67"
$t.Cell(7, 1).Range.Text = "11+45=56"
$t.Cell(7, 2).Range.Text = "21+40=61"
$t.Cell(7, 3).Range.Text = "89-12=77"
$t.Cell(7, 4).Range.Text = "64-54=10"
$t.Cell(7, 5).Range.Text = "16+39=55"
$t.Cell(8, 1).Range.Text = "72-55=17"
$t.Cell(8, 2).Range.Text = "8+62=70"
$t.Cell(8, 3).Range.Text = "85-75=10"
$t.Cell(8, 4).Range.Text = "77-45=32"
$t.Cell(8, 5).Range.Text = "4+35=39"
$t.Cell(9, 1).Range.Text = "52+40=92"
$t.Cell(9, 2).Range.Text = "10+1=11"
$t.Cell(9, 3).Range.Text = "74+13=87"
$t.Cell(9, 4).Range.Text = "99-52=47"
$t.Cell(9, 5).Range.Text = "17-4=13"
$t.Cell(10, 1).Range.Text = "94-57=37"
$t.Cell(10, 2).Range.Text = "1+22=23"
$t.Cell(10, 3).Range.Text = "97-23=74"
$t.Cell(10, 4).Range.Text = "44+20=64"
$t.Cell(10, 5).Range.Text = "51+48=99"
$t.Cell(11, 1).Range.Text = "54+24=78"
$t.Cell(11, 2).Range.Text = "89-74=15"
$t.Cell(11, 3).Range.Text = "12+34=46"
$t.Cell(11, 4).Range.Text = "99-40=59"
$t.Cell(11, 5).Range.Text = "26-9=17"
$t.Cell(12, 1).Range.Text = "48+14=62"
$t.Cell(12, 2).Range.Text = "28+71=99"
$t.Cell(12, 3).Range.Text = "74-37=37"
$t.Cell(12, 4).Range.Text = "64+6=70"
$t.Cell(12, 5).Range.Text = "6+71=77"
$t.Cell(13, 1).Range.Text = "75-73=2"
$t.Cell(13, 2).Range.Text = "45-17=28"
$t.Cell(13, 3).Range.Text = "56+13=69"
$t.Cell(13, 4).Range.Text = "44+39=83"
$t.Cell(13, 5).Range.Text = "65-45=20"
$t.Cell(14, 1).Range.Text = "11+14=25"
$t.Cell(14, 2).Range.Text = "85+2=87"
$t.Cell(14, 3).Range.Text = "84-27=57"
$t.Cell(14, 4).Range.Text = "42+1=43"
$t.Cell(14, 5).Range.Text = "33+1=34"
$t.Cell(15, 1).Range.Text = "78-27=51"
$t.Cell(15, 2).Range.Text = "18+19=37"
$t.Cell(15, 3).Range.Text = "82-36=46"
$t.Cell(15, 4).Range.Text = "92-92=0"
$t.Cell(15, 5).Range.Text = "77-16=61"
$t.Cell(16, 1).Range.Text = "3+31=34"
$t.Cell(16, 2).Range.Text = "96-36=60"
$t.Cell(16, 3).Range.Text = "83+4=87"
$t.Cell(16, 4).Range.Text = "60+12=72"
$t.Cell(16, 5).Range.Text = "84-33=51"
$t.Cell(17, 1).Range.Text = "8+84=92"
$t.Cell(17, 2).Range.Text = "21+59=80"
$t.Cell(17, 3).Range.Text = "49-4=45"
$t.Cell(17, 4).Range.Text = "87-39=48"
$t.Cell(17, 5).Range.Text = "69-6=63"
$t.Cell(18, 1).Range.Text = "32+49=81"
$t.Cell(18, 2).Range.Text = "45-25=20"
$t.Cell(18, 3).Range.Text = "22+26=48"
$t.Cell(18, 4).Range.Text = "42+42=84"
$t.Cell(18, 5).Range.Text = "89-81=8"
$t.Cell(19, 1).Range.Text = "42-29=13"
$t.Cell(19, 2).Range.Text = "79-25=54"
$t.Cell(19, 3).Range.Text = "60-19=41"
$t.Cell(19, 4).Range.Text = "41-31=10"
$t.Cell(19, 5).Range.Text = "89-15=74"
$t.Cell(20, 1).Range.Text = "76+16=92"
$t.Cell(20, 2).Range.Text = "49+45=94"
$t.Cell(20, 3).Range.Text = "44+17=61"
$t.Cell(20, 4).Range.Text = "72+4=76"
$t.Cell(20, 5).Range.Text = "92-19=73"
